# Apply cryptos list update (prices / 1h volume %) per commit
# "Updated cryptos list on Thu Feb 15 17:45:05 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '52.092.98'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '2.828.80'
$ws.Range('E3').Value = '  +3.19%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = "'355.68"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.55%  '
$ws.Range('D6').Value = "'113.50"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.25%  '
$ws.Range('D7').Value = "'0.547"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.39%  '
$ws.Range('D9').Value = "'0.607"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.85%  '
$ws.Range('D10').Value = "'42.05"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.71%  '
$ws.Range('D11').Value = "'0.0849"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.19%  '
$ws.Range('D12').Value = "'20.13"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.10%  '
$ws.Range('E13').Value = '  +1.30%  '
$ws.Range('D14').Value = "'7.77"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.16%  '
$ws.Range('D15').Value = '3.260.36'
$ws.Range('E15').Value = '  +2.80%  '
$ws.Range('D16').Value = '2.832.44'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('E17').Value = '  +2.47%  '
$ws.Range('D18').Value = '52.091.85'
$ws.Range('E18').Value = '  +1.14%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = "'7.29"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.94%  '
$ws.Range('D20').Value = "'3.16"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').Value = "'13.77"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.56%  '
$ws.Range('D22').Value = '0.0₃0994'
$ws.Range('E22').Value = '  +2.55%  '
$ws.Range('D23').Value = "'270.26"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.13%  '
$ws.Range('D24').Value = "'69.69"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.79%  '
$ws.Range('E25').Value = '  +5.59%  '
$ws.Range('D26').Value = "'26.69"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').Value = "'10.27"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.14%  '
$ws.Range('E29').Value = '  +1.54%  '
$ws.Range('E30').Value = '  +0.64%  '
$ws.Range('D31').Value = "'50.73"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.64%  '
$ws.Range('D32').Value = "'33.87"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('D33').Value = "'5.88"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +6.29%  '
$ws.Range('D34').Value = "'0.0444"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +28.95%  '
$ws.Range('D35').Value = "'0.0831"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('E36').Value = '  +0.11%  '
$ws.Range('D37').Value = "'2.09"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.39%  '
$ws.Range('D38').Value = "'4.88"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.38%  '
$ws.Range('D39').Value = "'18.43"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.50%  '
$ws.Range('D40').Value = "'3.21"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.34%  '
$ws.Range('E41').Value = '  +9.11%  '
$ws.Range('B42').Value = 'Monero'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D42').Value = "'128.00"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.37%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').Value = "'23.52"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.55%  '
$ws.Range('E44').Value = '  +1.87%  '
$ws.Range('E45').Value = '  +0.66%  '
$ws.Range('D46').Value = "'3.35"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.38%  '
$ws.Range('D47').Value = '2.044.41'
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('D49').Value = "'0.960"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +10.82%  '
$ws.Range('E50').Value = '  +3.41%  '
$ws.Range('D51').Value = "'60.54"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.61%  '
